# Auto-generated Excel COM-interop script to apply crypto price/volume updates
# matching the commit "Updated cryptos list on Thu Oct  5 18:30:13 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.468.86'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '''1.615.92'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.56%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '''210.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  -1.76%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '''22.71'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.24%  '
$ws.Range('E9').Value = '  +1.99%  '
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('D11').Value = '''0.0885'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').Value = '''1.847.06'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.45%  '
$ws.Range('D13').Value = '''1.613.23'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.43%  '
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('D15').Value = '''0.549'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.61%  '
$ws.Range('D16').Value = '''64.72'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '''27.500.34'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('D18').Value = '''229.59'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('D20').Value = '''7.51'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('D23').Value = '''10.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.40%  '
$ws.Range('E24').Value = '  +7.48%  '
$ws.Range('D25').Value = '''148.98'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').Value = '''0.111'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('B27').Value = 'BinanceUSD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D27').Value = '''1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.21%  '
$ws.Range('D28').Value = '''6.79'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.61%  '
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('D31').Value = '''0.0481'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('D32').Value = '''3.25'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.41%  '
$ws.Range('D33').Value = '''1.436.81'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('D34').Value = '''3.05'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.72%  '
$ws.Range('E35').Value = '  -3.61%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').Value = '''0.934'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.47%  '
$ws.Range('D38').Value = '''0.559'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.45%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').Value = '''0.858'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.09%  '
$ws.Range('D41').Value = '''69.11'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.01%  '
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('E43').Value = '  -2.55%  '
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('E45').Value = '  -2.35%  '
$ws.Range('E46').Value = '  -2.28%  '
$ws.Range('D47').Value = '''1.756.50'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.51%  '
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('D49').Value = '''86.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('D51').Value = '''0.0993'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.74%  '
